$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.687064829607016
$ws.Range("D2").Value = 4.999999999999999
$ws.Range("F2").Value = 2.999999999999999
$ws.Range("H2").Value = 1.521007298651291
$ws.Range("J2").Value = 0.7680720467421716
$ws.Range("L2").Value = 3.146244802641622

$ws.Range("B4").Value = 4.87629019523887
$ws.Range("D4").Value = 7.605036493256457
$ws.Range("F4").Value = 6.387017402085046
$ws.Range("H4").Value = 1.602928147569862
$ws.Range("J4").Value = 3.910178738044117
